# Small changes to performance profile plots.
#
# The K547:K555 "performance profile" summary formulas on each of the four
# benders_results_* sheets statistically summarise columns E:H of the
# sheet's data. They used to start their ranges at row 20 (skipping the
# first 18 data rows) - fix them to start at row 2, i.e. right after the
# header row, so the whole data set is included. The Q-column cells (only
# present on benders_results_0) already just AVERAGE(...) the K-column
# cells across all four sheets via a 3-D reference, so they recompute on
# their own once the K-column formulas are corrected - no need to touch
# them directly.
#
# Also nudge each sheet's scroll position / selected cell to match where
# the author's cursor ended up after making the change.

$wb = $excel.ActiveWorkbook

$sheetNames = @("benders_results_0", "benders_results_5", "benders_results_10", "benders_results_15")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("K547").Formula = "=100*COUNTIF(E2:E555,0)/ROWS(E2:E555)"
    $ws.Range("K548").Formula = "=COUNTIF(E2:E555,0)"
    $ws.Range("K549").Formula = "=AVERAGEIF(E2:E555,""<>*nan"")"
    $ws.Range("K550").Formula = "=AVERAGEIFS(E2:E555,E2:E555,""<>*nan"",E2:E555,""<>0"")"
    $ws.Range("K551").Formula = "=AVERAGEIF(F2:F555,""<>0"")"
    $ws.Range("K552").Formula = "=MEDIAN(IF(F2:F555<>0,F2:F555))"
    $ws.Range("K553").Formula = "=AVERAGEIF(G2:G555,""<>*inf"")"
    $ws.Range("K554").Formula = "=AVERAGE(H2:H555)"
    $ws.Range("K555").Formula = "=AVERAGEIF(H2:H555,""<7200"")"
}

# --- benders_results_0: topLeftCell K520 -> B520, selection Q552 -> K547 ---
$ws1 = $wb.Worksheets.Item("benders_results_0")
$ws1.Activate()
$excel.ActiveWindow.ScrollRow = 520
$excel.ActiveWindow.ScrollColumn = 2
$ws1.Range("K547").Select()

# --- benders_results_5: topLeftCell H536 -> H544, selection P547 -> I585 ---
$ws2 = $wb.Worksheets.Item("benders_results_5")
$ws2.Activate()
$excel.ActiveWindow.ScrollRow = 544
$excel.ActiveWindow.ScrollColumn = 8
$ws2.Range("I585").Select()

# --- benders_results_10: topLeftCell A523 (unchanged), selection J547 -> K547 ---
$ws3 = $wb.Worksheets.Item("benders_results_10")
$ws3.Activate()
$excel.ActiveWindow.ScrollRow = 523
$excel.ActiveWindow.ScrollColumn = 1
$ws3.Range("K547").Select()

# --- benders_results_15: topLeftCell A531 (unchanged), selection J550 -> K547 ---
$ws4 = $wb.Worksheets.Item("benders_results_15")
$ws4.Activate()
$excel.ActiveWindow.ScrollRow = 531
$excel.ActiveWindow.ScrollColumn = 1
$ws4.Range("K547").Select()

# benders_results_0 is the tab that was selected/active before and after the edit.
$ws1.Activate()
